# Add a new test step ("Input valid Value in the Account Number fields") to
# TestScenario_1 (New Account) in the Table1 list object.
#
# Before: row 4 (Step 3) held a generic placeholder "Valid value for required
# field" precondition/step/result. This edit repurposes row 4's UserAction /
# ExpectedResult text to be specific to the "Account Number" field, then
# inserts a brand-new row directly below (new row 5) that keeps the generic
# placeholder text as "Step 4" (reusing the existing shared strings), pushing
# the old "Step 4 / Save / Account created" row down to become row 6 (now
# relabeled "Step 5"), and every following row shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert one blank worksheet row above row 6, i.e. below the current row 4,
# which places a new empty row 5; shifts the old rows 5-16 down to 6-17.
$ws.Range("A5:J5").Insert()

# The table range doesn't auto-grow from a plain worksheet row insert here,
# so resize it explicitly to include the new row.
$tbl.Resize($ws.Range("A1:J17"))

# Row 4 (Step 3): replace the generic "Input valid value in the field" /
# "User should be able..." text with the Account-Number-specific wording.
$ws.Range("G4").Value = "Input valid Value in the Account Number fields"
$ws.Range("H4").Value = "User should be able to input value for the Account Number field."

# Row 5 (new): fill it back in with the generic placeholder text that used
# to live in row 4 (same as row 3's wording), now labelled "Step 4".
$ws.Range("E5").Value = "Valid value for required field  "
$ws.Range("F5").Value = "Step 4"
$ws.Range("G5").Value = "Input valid value in the   field."
$ws.Range("H5").Value = "User should be able to input value for the  field."

# Row 6 (shifted down from the old row 5): re-label its step number so the
# step sequence stays consistent after inserting the new Step 4.
$ws.Range("F6").Value = "Step 5"

# Match the author's final selection/scroll position on the sheet.
$ws.Activate()
$ws.Range("F6").Select()
